# Finish ethics checklist (#14)
# Answer "yes" to the row-9 question ("Does it involve tracking or
# observation of participants? ...") and bring the sheet view in line with
# where a user would naturally land after typing that answer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = "Yes"

# Reflect the scrolled/selected state left behind by the edit.
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("B10").Select()
